$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new C (nombre_aides) and new E (montant_total) values.
# D (nombre_entreprises) is unchanged for every row in this update.
$updates = @(
    @{Row=2;   C=766326;  E=1429219098}
    @{Row=3;   C=791;     E=2233378}
    @{Row=48;  C=150634;  E=275739556}
    @{Row=112; C=145230;  E=716365950}
    @{Row=121; C=1306266; E=2275067309}
    @{Row=129; C=633607;  E=3431371573}
    @{Row=131; C=378;     E=19427930}
    @{Row=132; C=585849;  E=3467997594}
    @{Row=136; C=26694;   E=144319346}
    @{Row=137; C=51;      E=2267833}
    @{Row=154; C=18462;   E=73574082}
    @{Row=186; C=236828;  E=1189788011}
    @{Row=221; C=135499;  E=681875940}
    @{Row=240; C=205914;  E=1069370264}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("C$r").Value = $u.C
    $ws.Range("E$r").Value = $u.E
}

$wb.Save()
